$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on target cells so values keep their exact original string representation
# (these cells were stored as inline/shared text, e.g. "40.60", not numbers)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"

# Apply updated Price (D) / Volume(1h) (E) values from the latest coinranking.com snapshot
$ws.Range("D2").Value = "329.34"
$ws.Range("E2").Value = "6.84%"
$ws.Range("D3").Value = "40.43"
$ws.Range("E3").Value = "12.59%"
$ws.Range("D4").Value = "5.940"
$ws.Range("E4").Value = "16.22%"
$ws.Range("D5").Value = "0.08143"
$ws.Range("E5").Value = "5.53%"
$ws.Range("D6").Value = "4.555"
$ws.Range("E6").Value = "3.87%"
$ws.Range("D7").Value = "8.751"
$ws.Range("E7").Value = "5.58%"
$ws.Range("D8").Value = "1.950"
$ws.Range("E8").Value = "5.18%"
$ws.Range("D9").Value = "2.999"
$ws.Range("E9").Value = "1.95%"
$ws.Range("D10").Value = "0.9436"
$ws.Range("E10").Value = "2.54%"
$ws.Range("D11").Value = "0.1311"
$ws.Range("E11").Value = "16.08%"
$ws.Range("D12").Value = "0.1992"
$ws.Range("E12").Value = "7.49%"
$ws.Range("D13").Value = "0.09232"
$ws.Range("E13").Value = "5.43%"
$ws.Range("E14").Value = "3.41%"
$ws.Range("D15").Value = "0.09618"
$ws.Range("E15").Value = "0.87%"
$ws.Range("D16").Value = "0.001318"
$ws.Range("E16").Value = "-4.55%"
$ws.Range("D17").Value = "0.006046"
$ws.Range("E17").Value = "-0.98%"
$ws.Range("D18").Value = "3.376"
$ws.Range("E18").Value = "0.39%"
$ws.Range("D19").Value = "0.3498"
$ws.Range("E19").Value = "1.53%"
$ws.Range("D20").Value = "7.729"
$ws.Range("E20").Value = "22.39%"
$ws.Range("D21").Value = "0.1437"
$ws.Range("E21").Value = "10.66%"
$ws.Range("D22").Value = "0.2445"
$ws.Range("E22").Value = "5.75%"
$ws.Range("D23").Value = "0.04424"
$ws.Range("E23").Value = "2.03%"
$ws.Range("E24").Value = "4.24%"
$ws.Range("D25").Value = "0.004373"
$ws.Range("E25").Value = "2.62%"
$ws.Range("D26").Value = "0.0001190"
$ws.Range("E26").Value = "-10.70%"
$ws.Range("D27").Value = "0.0003983"
$ws.Range("E27").Value = "37.27%"
$ws.Range("E39").Value = "19.03%"
$ws.Range("E40").Value = "7.90%"
$ws.Range("D41").Value = "0.007600"
$ws.Range("E41").Value = "0.81%"
$ws.Range("D42").Value = "0.1432"
$ws.Range("E42").Value = "6.15%"
$ws.Range("D43").Value = "0.008868"
$ws.Range("E43").Value = "4.11%"
$ws.Range("D44").Value = "0.002060"
$ws.Range("E44").Value = "-0.67%"
$ws.Range("D45").Value = "0.01049"
$ws.Range("E45").Value = "24.73%"
$ws.Range("D46").Value = "0.00006847"
$ws.Range("E46").Value = "5.96%"
$ws.Range("D47").Value = "0.00000000749"
$ws.Range("E47").Value = "-0.23%"
$ws.Range("D48").Value = "0.002893"
$ws.Range("E48").Value = "-12.37%"
$ws.Range("D49").Value = "0.001797"
$ws.Range("E49").Value = "24.47%"
$ws.Range("D50").Value = "0.00002097"
$ws.Range("E50").Value = "-0.23%"
$ws.Range("D51").Value = "0.0001997"
$ws.Range("E51").Value = "-0.23%"
